# ---------------------------------------------------------------------------
# 1. Add the new FAQ answer paragraph right after the "Not right now..."
#    paragraph.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "Not right now. I have a great community of friends here, and I’d rather not leave!",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Not right now...' paragraph."
}

# Insert a brand-new paragraph right after the found one, then fill it in.
$rng.InsertParagraphAfter()
$newPara = $d.Range($rng.End + 1, $rng.End + 1)
$newPara.Text = "I’m willing to commute to any location that is accessible via some combination of Metro + Bike + Scooter."

# ---------------------------------------------------------------------------
# 2. Make the "Heading 1" style's bottom border more prominent (swap the pale
#    gray dotted rule for the themed accent-1 dotted rule), and bump the
#    rsid stamps on the Heading1 / Heading1Char style definitions, exactly
#    like Word does when a style is touched in a new editing session.
# ---------------------------------------------------------------------------
$xml = $d.WordOpenXML

$oldHeading1 = '<w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00ED2422"/><w:pPr><w:keepNext/><w:keepLines/><w:pBdr><w:bottom w:val="dotted" w:sz="4" w:space="1" w:color="E7E6E6" w:themeColor="background2"/></w:pBdr>'
$newHeading1 = '<w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="007B15D0"/><w:pPr><w:keepNext/><w:keepLines/><w:pBdr><w:bottom w:val="dotted" w:sz="4" w:space="1" w:color="8EAADB" w:themeColor="accent1" w:themeTint="99"/></w:pBdr>'

if ($xml.IndexOf($oldHeading1) -lt 0) {
    throw "Heading1 style block not found in WordOpenXML."
}
$xml = $xml.Replace($oldHeading1, $newHeading1)

$oldHeading1Char = '<w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="00ED2422"/><w:rPr>'
$newHeading1Char = '<w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="007B15D0"/><w:rPr>'

if ($xml.IndexOf($oldHeading1Char) -lt 0) {
    throw "Heading1Char style block not found in WordOpenXML."
}
$xml = $xml.Replace($oldHeading1Char, $newHeading1Char)

$d.WordOpenXML = $xml

Write-Output "edit applied"
